$wb = $excel.ActiveWorkbook

# --- "URL" sheet: remove the JCVZ-997 entry (row 11). The row below
# (JCVZ-998) shifts up to become the new row 11. ---
$wsUrl = $wb.Worksheets.Item("URL")
$wsUrl.Rows.Item(11).Delete() | Out-Null

# --- "Result" sheet: remove the corresponding JCVZ-997 row (row 11);
# the JCVZ-998 row shifts up to row 11. Also remove the whole "Due Date"
# column (N), which is being dropped from the report. ---
$wsResult = $wb.Worksheets.Item("Result")
$wsResult.Rows.Item(11).Delete() | Out-Null
$wsResult.Columns.Item(14).Delete() | Out-Null

# --- restore the active selection on the URL sheet ---
$wsUrl.Range("H14").Select() | Out-Null
